# Auto update Excel log 2026-02-04 14:11:25
# Appends newly-logged sensor readings to the PIR, Humidity and Temperature
# sheets, plus the corresponding door open/close alerts on the ALERTS sheet.
#
# Column layout on every sheet: A Date | B Timestamp | C Hour | D Location | E Value | F Status
# Column A holds a "YYYY-MM-DD" string; Excel would otherwise reinterpret it
# as a real date serial, so it is written with a leading quote (forcing text)
# and the formatting that the quote leaves behind is cleared straight after.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# PIR sheet: rows 113-125 (No Motion / Inactive)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("PIR")

$pirTimes = @(
    "14:10:20","14:10:21","14:10:26","14:10:31","14:10:36","14:10:41",
    "14:10:46","14:10:51","14:10:56","14:11:01","14:11:06","14:11:11","14:11:16"
)

$rowCount = $pirTimes.Length
$startRow = 113
$endRow = $startRow + $rowCount - 1
$rng = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 6))

$out = New-Object 'object[,]' $rowCount, 6
for ($r = 0; $r -lt $rowCount; $r++) {
    $out[$r, 0] = "'2026-02-04"
    $out[$r, 1] = $pirTimes[$r]
    $out[$r, 2] = "14:00"
    $out[$r, 3] = "Bathroom"
    $out[$r, 4] = "No Motion"
    $out[$r, 5] = "Inactive"
}

$rng.Value = $out
$rng.ClearFormats()

# ---------------------------------------------------------------------------
# Humidity sheet: rows 89-100 (Active, with %RH readings)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Humidity")

$humTimes = @(
    "14:10:19","14:10:20","14:10:30","14:10:35","14:10:40","14:10:45",
    "14:10:50","14:10:55","14:11:00","14:11:05","14:11:10","14:11:15"
)
$humValues = @(
    "77.8%","76.7%","76.7%","77.6%","76.7%","77.8%",
    "76.9%","77.6%","76.6%","77.6%","76.5%","77.5%"
)

$rowCount = $humTimes.Length
$startRow = 89
$endRow = $startRow + $rowCount - 1
$rng = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 6))

$out = New-Object 'object[,]' $rowCount, 6
for ($r = 0; $r -lt $rowCount; $r++) {
    $out[$r, 0] = "'2026-02-04"
    $out[$r, 1] = $humTimes[$r]
    $out[$r, 2] = "14:00"
    $out[$r, 3] = "Bathroom"
    $out[$r, 4] = "'" + $humValues[$r]
    $out[$r, 5] = "Active"
}

$rng.Value = $out
$rng.ClearFormats()

# ---------------------------------------------------------------------------
# Temperature sheet: rows 89-100 (Active, with Celsius readings)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Temperature")

$tempTimes = @(
    "14:10:20","14:10:21","14:10:30","14:10:35","14:10:40","14:10:45",
    "14:10:51","14:10:55","14:11:01","14:11:05","14:11:11","14:11:15"
)
$tempValues = @(
    "24.8C","24.8C","24.8C","24.8C","24.8C","24.8C",
    "24.9C","24.8C","24.8C","24.9C","24.9C","24.9C"
)

$rowCount = $tempTimes.Length
$startRow = 89
$endRow = $startRow + $rowCount - 1
$rng = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 6))

$out = New-Object 'object[,]' $rowCount, 6
for ($r = 0; $r -lt $rowCount; $r++) {
    $out[$r, 0] = "'2026-02-04"
    $out[$r, 1] = $tempTimes[$r]
    $out[$r, 2] = "14:00"
    $out[$r, 3] = "Bathroom"
    $out[$r, 4] = $tempValues[$r]
    $out[$r, 5] = "Active"
}

$rng.Value = $out
$rng.ClearFormats()

# ---------------------------------------------------------------------------
# ALERTS sheet: rows 2-4 (Bathroom Door ENTER/EXIT events)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALERTS")

$alertTimes  = @("14:11:12","14:11:17","14:11:19")
$alertValues = @("ENTER","EXIT","ENTER")
$alertStatus = @("User ENTERED Bathroom","User EXITED Bathroom","User ENTERED Bathroom")

$rowCount = $alertTimes.Length
$startRow = 2
$endRow = $startRow + $rowCount - 1
$rng = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 6))

$out = New-Object 'object[,]' $rowCount, 6
for ($r = 0; $r -lt $rowCount; $r++) {
    $out[$r, 0] = "'2026-02-04"
    $out[$r, 1] = $alertTimes[$r]
    $out[$r, 2] = "14:00"
    $out[$r, 3] = "Bathroom Door"
    $out[$r, 4] = $alertValues[$r]
    $out[$r, 5] = $alertStatus[$r]
}

$rng.Value = $out
$rng.ClearFormats()
